# Weekly update: insert a new week's Jengibre price row (Vega Central
# Mapocho de Santiago) ahead of the existing history, pushing every
# older row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 66, shifting rows 66:108 down to 67:109.
$ws.Rows.Item(66).Insert()

# Populate the new row 66 with this week's record.
$ws.Range("A66").Value = 9
$ws.Range("B66").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C66").Value = "Metropolitana"
$ws.Range("D66").Value = 44830
$ws.Range("E66").Value = 13
$ws.Range("F66").Value = 100114007
$ws.Range("G66").Value = "Jengibre"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 450
$ws.Range("K66").Value = 14000
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = 14556
$ws.Range("N66").Value = "`$/caja 13 kilos"
$ws.Range("O66").Value = "Perú"
$ws.Range("P66").Value = 1120
$ws.Range("Q66").Value = 13
$ws.Range("R66").Value = "Hortaliza"
